$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on Price cells so numeric-looking strings
# (e.g. "42.604.21", "309.56") are stored as text, matching the source data.
$priceCells = @('D2', 'D3', 'D5', 'D6', 'D7', 'D9', 'D10', 'D11', 'D12', 'D14', 'D15', 'D16', 'D17', 'D18', 'D19', 'D20', 'D21', 'D22', 'D23', 'D24', 'D25', 'D26', 'D28', 'D29', 'D30', 'D31', 'D32', 'D33', 'D34', 'D35', 'D36', 'D39', 'D40', 'D41', 'D42', 'D44', 'D45', 'D46', 'D47', 'D48', 'D49', 'D50', 'D51')
foreach ($c in $priceCells) { $ws.Range($c).NumberFormat = "@" }

# Apply the updated values
$ws.Range('D2').Value = '42.604.21'
$ws.Range('E2').Value = '  -0.29%  '
$ws.Range('D3').Value = '2.518.37'
$ws.Range('E3').Value = '  -1.55%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '309.56'
$ws.Range('E5').Value = '  +2.28%  '
$ws.Range('D6').Value = '96.61'
$ws.Range('E6').Value = '  -0.62%  '
$ws.Range('D7').Value = '0.587'
$ws.Range('E7').Value = '  +2.12%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').Value = '0.540'
$ws.Range('E9').Value = '  -1.21%  '
$ws.Range('D10').Value = '36.75'
$ws.Range('E10').Value = '  +0.52%  '
$ws.Range('D11').Value = '0.0814'
$ws.Range('E11').Value = '  +0.75%  '
$ws.Range('D12').Value = '7.76'
$ws.Range('E12').Value = '  +2.92%  '
$ws.Range('E13').Value = '  -3.53%  '
$ws.Range('D14').Value = '2.902.48'
$ws.Range('E14').Value = '  -1.79%  '
$ws.Range('D15').Value = '15.84'
$ws.Range('E15').Value = '  +8.59%  '
$ws.Range('D16').Value = '2.506.67'
$ws.Range('E16').Value = '  -2.86%  '
$ws.Range('D17').Value = '0.859'
$ws.Range('E17').Value = '  -2.63%  '
$ws.Range('D18').Value = '42.565.96'
$ws.Range('E18').Value = '  -0.52%  '
$ws.Range('D19').Value = '12.99'
$ws.Range('E19').Value = '  -4.60%  '
$ws.Range('D20').Value = '0.0₃0974'
$ws.Range('E20').Value = '  -1.45%  '
$ws.Range('D21').Value = '6.47'
$ws.Range('E21').Value = '  -2.17%  '
$ws.Range('D22').Value = '71.57'
$ws.Range('E22').Value = '  -0.06%  '
$ws.Range('D23').Value = '253.95'
$ws.Range('E23').Value = '  -1.08%  '
$ws.Range('D24').Value = '2.94'
$ws.Range('E24').Value = '  +0.02%  '
$ws.Range('D25').Value = '2.05'
$ws.Range('E25').Value = '  -2.16%  '
$ws.Range('D26').Value = '27.15'
$ws.Range('E26').Value = '  -3.46%  '
$ws.Range('E27').Value = '  +0.08%  '
$ws.Range('D28').Value = '2.34'
$ws.Range('E28').Value = '  +10.94%  '
$ws.Range('D29').Value = '10.20'
$ws.Range('E29').Value = '  +1.20%  '
$ws.Range('D30').Value = '37.65'
$ws.Range('E30').Value = '  -4.12%  '
$ws.Range('D31').Value = '5.96'
$ws.Range('E31').Value = '  -0.71%  '
$ws.Range('D32').Value = '154.20'
$ws.Range('E32').Value = '  -1.07%  '
$ws.Range('D33').Value = '19.20'
$ws.Range('E33').Value = '  +5.47%  '
$ws.Range('D34').Value = '3.29'
$ws.Range('E34').Value = '  -1.59%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').Value = '2.08'
$ws.Range('E35').Value = '  -3.98%  '
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').Value = '0.0790'
$ws.Range('E36').Value = '  -1.60%  '
$ws.Range('E37').Value = '  -4.90%  '
$ws.Range('E38').Value = '  -1.35%  '
$ws.Range('B39').Value = 'EnergySwap'
$ws.Range('C39').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D39').Value = '24.40'
$ws.Range('E39').Value = '  -10.00%  '
$ws.Range('B40').Value = 'Stellar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D40').Value = '0.120'
$ws.Range('E40').Value = '  +0.47%  '
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D41').Value = '3.88'
$ws.Range('E41').Value = '  +0.74%  '
$ws.Range('B42').Value = 'NEARProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D42').Value = '3.40'
$ws.Range('E42').Value = '  +0.68%  '
$ws.Range('E43').Value = '  -0.20%  '
$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D44').Value = '0.0302'
$ws.Range('E44').Value = '  -1.07%  '
$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D45').Value = '0.998'
$ws.Range('E45').Value = '  -0.07%  '
$ws.Range('D46').Value = '2.038.19'
$ws.Range('E46').Value = '  -1.07%  '
$ws.Range('D47').Value = '84.81'
$ws.Range('E47').Value = '  -4.11%  '
$ws.Range('D48').Value = '8.99'
$ws.Range('E48').Value = '  -2.74%  '
$ws.Range('D49').Value = '2.756.67'
$ws.Range('E49').Value = '  -1.93%  '
$ws.Range('D50').Value = '73.07'
$ws.Range('E50').Value = '  -5.13%  '
$ws.Range('D51').Value = '0.191'
$ws.Range('E51').Value = '  +0.23%  '

# Restore default style on Price cells (keeps the text content, drops the
# temporary text format so the cell style matches the rest of the sheet).
foreach ($c in $priceCells) { $ws.Range($c).Style = "Normal" }
